# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newer scrape counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, sheet "展览"
$exhibitionUpdates = @{
    5  = 1172
    6  = 14375
    7  = 16684
    9  = 112
    10 = 21
    11 = 51
    20 = 39
    21 = 1272
    22 = 138
    23 = 71
    24 = 43
    25 = 26
    27 = 6798
    30 = 1127
    31 = 14
    33 = 5769
    37 = 4867
}

# Row -> new value for column F, sheet "全部类型"
$allTypesUpdates = @{
    5  = 1172
    6  = 14375
    7  = 16684
    9  = 112
    10 = 21
    11 = 51
    20 = 39
    21 = 1272
    22 = 138
    23 = 71
    25 = 43
    26 = 26
    28 = 6798
    31 = 1127
    32 = 14
    36 = 5769
    40 = 4867
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
